$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each changed cell value, matching the source diff exactly.
# Cells in column D that look like plain numbers need to be forced to
# Text format first so Excel stores them as strings (as in the original
# workbook) instead of auto-converting them to numeric values.

$ws.Range('D2').Value = '42.371.22'
$ws.Range('E2').Value = '  -0.89%  '

$ws.Range('D3').Value = '2.512.74'
$ws.Range('E3').Value = '  -2.29%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.996'
$ws.Range('E4').Value = '  -0.34%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '304.17'
$ws.Range('E5').Value = '  +0.72%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '96.83'

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.579'
$ws.Range('E7').Value = '  +0.77%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.16%  '

$ws.Range('E9').Value = '  -1.59%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.61'
$ws.Range('E10').Value = '  +0.72%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0810'
$ws.Range('E11').Value = '  -0.16%  '

$ws.Range('E12').Value = '  -0.80%  '

$ws.Range('E13').Value = '  -1.68%  '

$ws.Range('D14').Value = '2.889.73'
$ws.Range('E14').Value = '  -2.68%  '

$ws.Range('D15').Value = '2.523.22'
$ws.Range('E15').Value = '  -3.34%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '15.01'
$ws.Range('E16').Value = '  +4.92%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.860'
$ws.Range('E17').Value = '  -2.60%  '

$ws.Range('D18').Value = '42.571.09'
$ws.Range('E18').Value = '  -0.52%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.93'
$ws.Range('E19').Value = '  +0.32%  '

$ws.Range('D20').Value = '0.0₃0973'
$ws.Range('E20').Value = '  -2.46%  '

$ws.Range('E21').Value = '  -3.20%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '71.05'
$ws.Range('E22').Value = '  -1.20%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '250.76'
$ws.Range('E23').Value = '  -1.10%  '

$ws.Range('E24').Value = '  -1.20%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.04'
$ws.Range('E25').Value = '  -4.43%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '26.99'
$ws.Range('E26').Value = '  -6.66%  '

$ws.Range('E27').Value = '  +0.15%  '

$ws.Range('E28').Value = '  +10.79%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.32'
$ws.Range('E29').Value = '  +0.53%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '38.04'
$ws.Range('E30').Value = '  +0.78%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.96'
$ws.Range('E31').Value = '  -0.82%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '156.76'
$ws.Range('E32').Value = '  +1.03%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0791'
$ws.Range('E33').Value = '  -1.36%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.27'
$ws.Range('E34').Value = '  -4.06%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.07'
$ws.Range('E35').Value = '  -4.33%  '

$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.62'
$ws.Range('E36').Value = '  -4.58%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '18.41'
$ws.Range('E37').Value = '  +0.77%  '

$ws.Range('E38').Value = '  +1.61%  '

$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.119'
$ws.Range('E39').Value = '  -0.47%  '

$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '23.89'
$ws.Range('E40').Value = '  +4.18%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.08'
$ws.Range('E41').Value = '  -0.46%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.39'
$ws.Range('E42').Value = '  -0.84%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.84'
$ws.Range('E43').Value = '  -1.18%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.993'
$ws.Range('E44').Value = '  -0.53%  '

$ws.Range('E45').Value = '  -3.39%  '

$ws.Range('D46').Value = '2.032.52'
$ws.Range('E46').Value = '  -2.44%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '85.14'
$ws.Range('E47').Value = '  -0.20%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.95'
$ws.Range('E48').Value = '  -2.97%  '

$ws.Range('D49').Value = '2.753.04'
$ws.Range('E49').Value = '  -2.55%  '

$ws.Range('E50').Value = '  -0.99%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '101.49'
$ws.Range('E51').Value = '  -4.39%  '
